# Generate Report for Handoff
#
# The 281112d5-af5a-4b77-ac95-317023dcf2e9.md entry moved from "Handed
# back: in sync with en-US" to "Ready for handoff" (handoff re-run), while
# the 72dfbeb8-593a-400d-a8d4-77237649c29b.md entry's data simply shifted
# position (row swap) in every sheet. A stale-handback warning is also
# recorded in the Error Detail column of the locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
$ws.Range("B2").Value = "e2e\72dfbeb8-593a-400d-a8d4-77237649c29b.md"

$ws.Range("A3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$ws.Range("B3").Value = "e2e\281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-15 22:45:24"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\72dfbeb8-593a-400d-a8d4-77237649c29b.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\281112d5-af5a-4b77-ac95-317023dcf2e9.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
$ws.Range("G2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.c6f187c302caabd680fd62d118c073449a94aea0.zh-cn.xlf"
$ws.Range("I2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
$ws.Range("J2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.c6f187c302caabd680fd62d118c073449a94aea0.zh-cn.xlf"

$ws.Range("A3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.87c83b646adcf3a265a49c0c022b9aa4c78b8642.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-15 22:45:18"
$ws.Range("I3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$ws.Range("J3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.87c83b646adcf3a265a49c0c022b9aa4c78b8642.zh-cn.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2a78fabeb79efd4b344c8702d0df0ac6f1118e4/e2e/281112d5-af5a-4b77-ac95-317023dcf2e9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2dd27a90fa96dc63cd9b5c8e32c944ea27dc9bb/e2e/281112d5-af5a-4b77-ac95-317023dcf2e9.md."

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
    }
}

$ws.Range("P1").ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
$ws.Range("G2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.c6f187c302caabd680fd62d118c073449a94aea0.de-de.xlf"
$ws.Range("I2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
$ws.Range("J2").Value = "72dfbeb8-593a-400d-a8d4-77237649c29b.c6f187c302caabd680fd62d118c073449a94aea0.de-de.xlf"

$ws.Range("A3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.87c83b646adcf3a265a49c0c022b9aa4c78b8642.de-de.xlf"
$ws.Range("H3").Value = "2016-08-15 22:45:24"
$ws.Range("I3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$ws.Range("J3").Value = "281112d5-af5a-4b77-ac95-317023dcf2e9.87c83b646adcf3a265a49c0c022b9aa4c78b8642.de-de.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2a78fabeb79efd4b344c8702d0df0ac6f1118e4/e2e/281112d5-af5a-4b77-ac95-317023dcf2e9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2dd27a90fa96dc63cd9b5c8e32c944ea27dc9bb/e2e/281112d5-af5a-4b77-ac95-317023dcf2e9.md."

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
    }
}

$ws.Range("P1").ColumnWidth = 39.1666666666667
